# Apply edits described by the diff:
# 1. Rename the "students" sheet to "test"
# 2. Update header cell B1 from "University" to "Hello"
# 3. Update cell A2 from "Munerah Alzaidan" to "Munerah Alzaidan2"
# 4. Move the active selection to B1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "test"

# Update the student name in A2 (set first so it takes the earlier shared-string index)
$ws.Range("A2").Value = "Munerah Alzaidan2"

# Update the header value in B1
$ws.Range("B1").Value = "Hello"

# Move / update the selection to B1
$ws.Activate()
$ws.Range("B1").Select()
